$d = $word.ActiveDocument

# --- Create the three new character styles (wdStyleTypeCharacter = 2) ---

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "2022 Fechas de la campaña ..." run (4x) ---

$rng = $d.Content
$campaignText = "2022 Fechas de la campaña para constelación de perseo: 16-25 de enero, 7-16 de noviembre, 6-15 de diciembre"
while ($rng.Find.Execute($campaignText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# --- Apply GaNParagraph to the "Usted está participando ..." run ---

$rng = $d.Content
$paragraphText = "Usted está participando en una campaña mundial para observar y registrar las estrellas visibles más débiles como un medio para medir la contaminación lumínica en un lugar determinado. Localizando y observando la  constelación de perseo en el cielo nocturno y comparándolo con las cartas estelares, la gente de todo el mundo aprenderán cómo las luces de su comunidad contribuyen a la contaminación lumínica. Sus contribuciones a la base de datos en línea documentarán el cielo nocturno visible."
if ($rng.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "(http://amper.ped.muni.cz/..." run ---

$rng = $d.Content
$linkText = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
if ($rng.Find.Execute($linkText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNLinks"
}
